$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Erros - Tabela Geral" (row 11) as done: copy the visual format
# (fill/border) from an existing "ok" row and flip its status text.
$ws.Range("B4").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("B11").Value = "ok"

# Remove the completed "Correcao Transacao Cartao" row entirely; everything
# below shifts up by one row (table/autofilter/dimension follow along).
$ws.Rows(12).Delete()

# Restore the selection to where the author left off after the edit.
$ws.Range("G21").Select()
